# fix: fixed formatting when scrapping floating point numbers
#
# The "Importe" column (H) was scraped with Argentine-locale formatting
# (thousands separator "." and decimal separator ",", e.g. "36.550,00").
# Re-format every such value to a plain "plain" numeric-looking string
# with a dot decimal separator and no thousands separators
# (e.g. "36550.00"), while keeping the cell's content as TEXT (it must
# stay a literal string, not become a real Excel number).
#
# Separately, a handful of "Razon social" / "Nombre Fantasia" entries used
# a comma where a period was intended (multiple people/entities listed in
# one field); those commas are normalized to periods too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

# --- 1) Re-format Argentine-style decimal numbers stored as text -----------
# Pattern: optional thousands groups separated by '.', then ',' + 2 decimals
# e.g. "36.550,00", "530,00", "2.873.935,06", "9,00"
$numberPattern = '^\d{1,3}(\.\d{3})*,\d{2}$'

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val -and $val -match $numberPattern) {
            $newVal = ($val -replace '\.', '') -replace ',', '.'
            # Force the cell to stay TEXT: without this, assigning a
            # numeric-looking string turns the cell into a real number
            # (and drops the trailing zeros), which is not what we want.
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            # Reset back to the default style so we don't leave a
            # lingering custom number format applied to the cell.
            $cell.Style = "Normal"
        }
    }
}

# --- 2) Normalize stray commas to periods in specific name fields ----------
$nameFixes = @{
    "FERNANDEZ, MARIO HUGO" = "FERNANDEZ. MARIO HUGO"
    "MONROY, AGUSTIN ALEJANDRO" = "MONROY. AGUSTIN ALEJANDRO"
    "MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "RICCOTTI, MARIANA EDITH" = "RICCOTTI. MARIANA EDITH"
    "ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "MERCANZINI, GASTON ARIEL" = "MERCANZINI. GASTON ARIEL"
}

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val -and $nameFixes.ContainsKey($val)) {
            $cell.Value = $nameFixes[$val]
        }
    }
}
